$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data block from 3 rows (ECs/FAPs/sCs x ECs only) to the full
# 3x3 sending-cluster x target-cluster grid (rows 2-10), refreshed per Dr Hou advice.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serping1"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.90597433333333
$ws.Range("H2").Value = 35.717923
$ws.Range("I2").Value = 0.008895149679642379
$ws.Range("J2").Value = 0.008895149679642379
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.47090666666666
$ws.Range("N2").Value = 52.41271999999999
$ws.Range("O2").Value = 0.9803569739482672
$ws.Range("P2").Value = 0.9803569739482673
$ws.Range("Q2").Value = 208.0081663533955
$ws.Range("R2").Value = 1872.07349718056
$ws.Range("S2").Value = 0.0087204220227511
$ws.Range("T2").Value = 0.008720422022751102

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serping1"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.90597433333333
$ws.Range("H3").Value = 35.717923
$ws.Range("I3").Value = 0.008895149679642379
$ws.Range("J3").Value = 0.008895149679642379
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.16459
$ws.Range("N3").Value = 0.49377
$ws.Range("O3").Value = 0.0092357516081294
$ws.Range("P3").Value = 0.0092357516081294
$ws.Range("Q3").Value = 1.959604315523333
$ws.Range("R3").Value = 17.63643883971
$ws.Range("S3").Value = 0.00008215339295830881
$ws.Range("T3").Value = 0.00008215339295830881

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Serping1"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.90597433333333
$ws.Range("H4").Value = 35.717923
$ws.Range("I4").Value = 0.008895149679642379
$ws.Range("J4").Value = 0.008895149679642379
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1854676666666667
$ws.Range("N4").Value = 0.556403
$ws.Range("O4").Value = 0.01040727444360334
$ws.Range("P4").Value = 0.01040727444360334
$ws.Range("Q4").Value = 2.208173278996556
$ws.Range("R4").Value = 19.873559510969
$ws.Range("S4").Value = 0.00009257426393296859
$ws.Range("T4").Value = 0.00009257426393296859

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serping1"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1205.102620666667
$ws.Range("H5").Value = 3615.307862
$ws.Range("I5").Value = 0.9003520325209805
$ws.Range("J5").Value = 0.9003520325209804
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.47090666666666
$ws.Range("N5").Value = 52.41271999999999
$ws.Range("O5").Value = 0.9803569739482672
$ws.Range("P5").Value = 0.9803569739482673
$ws.Range("Q5").Value = 21054.23540942274
$ws.Range("R5").Value = 189488.1186848046
$ws.Range("S5").Value = 0.8826663940904402
$ws.Range("T5").Value = 0.8826663940904402

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Serping1"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1205.102620666667
$ws.Range("H6").Value = 3615.307862
$ws.Range("I6").Value = 0.9003520325209805
$ws.Range("J6").Value = 0.9003520325209804
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.16459
$ws.Range("N6").Value = 0.49377
$ws.Range("O6").Value = 0.0092357516081294
$ws.Range("P6").Value = 0.0092357516081294
$ws.Range("Q6").Value = 198.3478403355267
$ws.Range("R6").Value = 1785.13056301974
$ws.Range("S6").Value = 0.00831542773223822
$ws.Range("T6").Value = 0.008315427732238218

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Serping1"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1205.102620666667
$ws.Range("H7").Value = 3615.307862
$ws.Range("I7").Value = 0.9003520325209805
$ws.Range("J7").Value = 0.9003520325209804
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1854676666666667
$ws.Range("N7").Value = 0.556403
$ws.Range("O7").Value = 0.01040727444360334
$ws.Range("P7").Value = 0.01040727444360334
$ws.Range("Q7").Value = 223.5075711489318
$ws.Range("R7").Value = 2011.568140340386
$ws.Range("S7").Value = 0.009370210698301925
$ws.Range("T7").Value = 0.009370210698301925

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Serping1"
$ws.Range("C8").Value = "Selp"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 121.4707743333333
$ws.Range("H8").Value = 364.412323
$ws.Range("I8").Value = 0.09075281779937723
$ws.Range("J8").Value = 0.09075281779937722
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.47090666666666
$ws.Range("N8").Value = 52.41271999999999
$ws.Range("O8").Value = 0.9803569739482672
$ws.Range("P8").Value = 0.9803569739482673
$ws.Range("Q8").Value = 2122.204561105395
$ws.Range("R8").Value = 19099.84104994856
$ws.Range("S8").Value = 0.0889701578350759
$ws.Range("T8").Value = 0.08897015783507589

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Serping1"
$ws.Range("C9").Value = "Selp"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 121.4707743333333
$ws.Range("H9").Value = 364.412323
$ws.Range("I9").Value = 0.09075281779937723
$ws.Range("J9").Value = 0.09075281779937722
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.16459
$ws.Range("N9").Value = 0.49377
$ws.Range("O9").Value = 0.0092357516081294
$ws.Range("P9").Value = 0.0092357516081294
$ws.Range("Q9").Value = 19.99287474752333
$ws.Range("R9").Value = 179.93587272771
$ws.Range("S9").Value = 0.0008381704829328727
$ws.Range("T9").Value = 0.0008381704829328726

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Serping1"
$ws.Range("C10").Value = "Selp"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 121.4707743333333
$ws.Range("H10").Value = 364.412323
$ws.Range("I10").Value = 0.09075281779937723
$ws.Range("J10").Value = 0.09075281779937722
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1854676666666667
$ws.Range("N10").Value = 0.556403
$ws.Range("O10").Value = 0.01040727444360334
$ws.Range("P10").Value = 0.01040727444360334
$ws.Range("Q10").Value = 22.52890108379656
$ws.Range("R10").Value = 202.760109754169
$ws.Range("S10").Value = 0.0009444894813684492
$ws.Range("T10").Value = 0.0009444894813684491

